{"js": "// Replace the date line and every three-digit-by-one-digit multiplication\n// fact cell in the table with its updated value, per the commit diff.\nconst replacements = [\n  [\"2024-03-26 Tuesday\", \"2024-03-27 Wednesday\"],\n  [\"985\\u00d74=3940\", \"119\\u00d75=595\"],\n  [\"299\\u00d74=1196\", \"770\\u00d72=1540\"],\n  [\"369\\u00d77=2583\", \"942\\u00d75=4710\"],\n  [\"937\\u00d74=3748\", \"615\\u00d78=4920\"],\n  [\"831\\u00d77=5817\", \"252\\u00d74=1008\"],\n  [\"515\\u00d75=2575\", \"394\\u00d77=2758\"],\n  [\"402\\u00d79=3618\", \"978\\u00d72=1956\"],\n  [\"321\\u00d77=2247\", \"963\\u00d77=6741\"],\n  [\"965\\u00d73=2895\", \"204\\u00d75=1020\"],\n  [\"868\\u00d78=6944\", \"995\\u00d76=5970\"],\n  [\"916\\u00d77=6412\", \"573\\u00d75=2865\"],\n  [\"718\\u00d77=5026\", \"273\\u00d78=2184\"],\n  [\"110\\u00d73=330\", \"306\\u00d78=2448\"],\n  [\"437\\u00d79=3933\", \"509\\u00d77=3563\"],\n  [\"818\\u00d79=7362\", \"749\\u00d74=2996\"],\n  [\"721\\u00d73=2163\", \"497\\u00d75=2485\"],\n  [\"310\\u00d79=2790\", \"234\\u00d73=702\"],\n  [\"606\\u00d76=3636\", \"417\\u00d78=3336\"],\n  [\"466\\u00d76=2796\", \"162\\u00d73=486\"],\n  [\"202\\u00d73=606\", \"114\\u00d79=1026\"],\n  [\"575\\u00d76=3450\", \"845\\u00d77=5915\"],\n  [\"555\\u00d76=3330\", \"735\\u00d77=5145\"],\n  [\"668\\u00d79=6012\", \"703\\u00d78=5624\"],\n  [\"679\\u00d72=1358\", \"128\\u00d73=384\"],\n  [\"576\\u00d76=3456\", \"577\\u00d72=1154\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every three-digit-by-one-digit multiplication\n# fact cell in the table with its updated value, per the commit diff.\n$pairs = @(\n  @(\"2024-03-26 Tuesday\", \"2024-03-27 Wednesday\"),\n  @(\"985\u00d74=3940\", \"119\u00d75=595\"),\n  @(\"299\u00d74=1196\", \"770\u00d72=1540\"),\n  @(\"369\u00d77=2583\", \"942\u00d75=4710\"),\n  @(\"937\u00d74=3748\", \"615\u00d78=4920\"),\n  @(\"831\u00d77=5817\", \"252\u00d74=1008\"),\n  @(\"515\u00d75=2575\", \"394\u00d77=2758\"),\n  @(\"402\u00d79=3618\", \"978\u00d72=1956\"),\n  @(\"321\u00d77=2247\", \"963\u00d77=6741\"),\n  @(\"965\u00d73=2895\", \"204\u00d75=1020\"),\n  @(\"868\u00d78=6944\", \"995\u00d76=5970\"),\n  @(\"916\u00d77=6412\", \"573\u00d75=2865\"),\n  @(\"718\u00d77=5026\", \"273\u00d78=2184\"),\n  @(\"110\u00d73=330\", \"306\u00d78=2448\"),\n  @(\"437\u00d79=3933\", \"509\u00d77=3563\"),\n  @(\"818\u00d79=7362\", \"749\u00d74=2996\"),\n  @(\"721\u00d73=2163\", \"497\u00d75=2485\"),\n  @(\"310\u00d79=2790\", \"234\u00d73=702\"),\n  @(\"606\u00d76=3636\", \"417\u00d78=3336\"),\n  @(\"466\u00d76=2796\", \"162\u00d73=486\"),\n  @(\"202\u00d73=606\", \"114\u00d79=1026\"),\n  @(\"575\u00d76=3450\", \"845\u00d77=5915\"),\n  @(\"555\u00d76=3330\", \"735\u00d77=5145\"),\n  @(\"668\u00d79=6012\", \"703\u00d78=5624\"),\n  @(\"679\u00d72=1358\", \"128\u00d73=384\"),\n  @(\"576\u00d76=3456\", \"577\u00d72=1154\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $range = $d.Content\n  $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $found) {\n    Write-Output \"NOT FOUND: $oldText\"\n  }\n}\n"}
